$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.490.59"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.88%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.435.46"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.56%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "407.95"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.84"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.47%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.51%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.719"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +9.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.140"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +18.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.83"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.77%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.67"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.04"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.441.07"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.517.94"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.88%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.22%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000169"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +34.41%  "
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "Polygon"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.02"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.32%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "85.06"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +6.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "315.85"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +5.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.96"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.81%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "30.33"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +5.11%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.22"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.88"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "45.42"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +11.09%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +8.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.173"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.117"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.63"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.04%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0489"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.86"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.98"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.36"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.322"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +14.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "142.97"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.30%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.99"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.23%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.99"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.26"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.67"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.107.42"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.04"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +10.90%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.127"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +16.05%  "

Write-Host "Applied cryptos update"
